# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for specific leve rows across all eight job
# sheets, matching the values pulled from the latest price-data run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2831.6667
$ws.Range("I62").Value = 2598.8
$ws.Range("J62").Value = 2998
$ws.Range("K62").Value = 2598.8
$ws.Range("L62").Value = 2998
$ws.Range("M62").Value = -1974.8
$ws.Range("N62").Value = -4246

$ws.Range("H65").Value = 2831.6667
$ws.Range("I65").Value = 2598.8
$ws.Range("J65").Value = 2998
$ws.Range("K65").Value = 12994
$ws.Range("L65").Value = 14990
$ws.Range("M65").Value = -9874
$ws.Range("N65").Value = -21230

$ws.Range("H107").Value = 605.2857
$ws.Range("I107").Value = 373.75
$ws.Range("J107").Value = 914
$ws.Range("K107").Value = 373.75
$ws.Range("L107").Value = 914
$ws.Range("M107").Value = 1546.25
$ws.Range("N107").Value = -4754

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = 0

$ws.Range("H132").Value = 1269.9117
$ws.Range("I132").Value = 1054.2413
$ws.Range("K132").Value = 3162.7239
$ws.Range("M132").Value = -632.7239

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3916.61
$ws.Range("I32").Value = 3257.7778
$ws.Range("J32").Value = 6034.2856
$ws.Range("K32").Value = 3257.7778
$ws.Range("L32").Value = 6034.2856
$ws.Range("M32").Value = -2970.7778
$ws.Range("N32").Value = -6608.2856

$ws.Range("H45").Value = 1625.0667
$ws.Range("I45").Value = 1100.75
$ws.Range("J45").Value = 1815.7273
$ws.Range("K45").Value = 1100.75
$ws.Range("L45").Value = 1815.7273
$ws.Range("M45").Value = -723.75
$ws.Range("N45").Value = -2569.7273

$ws.Range("H74").Value = 1032.9565
$ws.Range("I74").Value = 872.625
$ws.Range("J74").Value = 2101.8333
$ws.Range("K74").Value = 872.625
$ws.Range("L74").Value = 2101.8333
$ws.Range("M74").Value = 1.375
$ws.Range("N74").Value = -3849.8333

$ws.Range("H77").Value = 1032.9565
$ws.Range("I77").Value = 872.625
$ws.Range("J77").Value = 2101.8333
$ws.Range("K77").Value = 4363.125
$ws.Range("L77").Value = 10509.1665
$ws.Range("M77").Value = 4.875
$ws.Range("N77").Value = -19245.1665

$ws.Range("H122").Value = 1658.6666
$ws.Range("I122").Value = 1616
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4848
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2398
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 2033
$ws.Range("I132").Value = 1649.5
$ws.Range("J132").Value = 2446
$ws.Range("K132").Value = 4948.5
$ws.Range("L132").Value = 7338
$ws.Range("M132").Value = -2418.5
$ws.Range("N132").Value = -12398

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1883.25
$ws.Range("I7").Value = 1177.6666
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 1177.6666
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -1064.6666
$ws.Range("N7").Value = -4226

$ws.Range("H86").Value = 102418.35
$ws.Range("I86").Value = 2255
$ws.Range("K86").Value = 2255
$ws.Range("M86").Value = -1132

$ws.Range("H89").Value = 102418.35
$ws.Range("I89").Value = 2255
$ws.Range("K89").Value = 11275
$ws.Range("M89").Value = -5659

$ws.Range("H94").Value = 406.59375
$ws.Range("I94").Value = 314.32144
$ws.Range("J94").Value = 1052.5
$ws.Range("K94").Value = 314.32144
$ws.Range("L94").Value = 1052.5
$ws.Range("M94").Value = 136.67856
$ws.Range("N94").Value = -1954.5

$ws.Range("H107").Value = 1393.3077
$ws.Range("I107").Value = 1264.579
$ws.Range("K107").Value = 1264.579
$ws.Range("M107").Value = 655.421

$ws.Range("H134").Value = 6875.154
$ws.Range("I134").Value = 7761.5454
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 23284.6362
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -20749.6362
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 846.1111
$ws.Range("I16").Value = 802.5714
$ws.Range("K16").Value = 802.5714
$ws.Range("M16").Value = -515.5714

$ws.Range("H31").Value = 1886.6957
$ws.Range("I31").Value = 1178
$ws.Range("J31").Value = 2264.6667
$ws.Range("K31").Value = 1178
$ws.Range("L31").Value = 2264.6667
$ws.Range("M31").Value = -883
$ws.Range("N31").Value = -2854.6667

$ws.Range("H34").Value = 1886.6957
$ws.Range("I34").Value = 1178
$ws.Range("J34").Value = 2264.6667
$ws.Range("K34").Value = 1178
$ws.Range("L34").Value = 2264.6667
$ws.Range("M34").Value = -976
$ws.Range("N34").Value = -2668.6667

$ws.Range("H99").Value = 1113625
$ws.Range("I99").Value = 5000505.5
$ws.Range("J99").Value = 3087.7144
$ws.Range("K99").Value = 5000505.5
$ws.Range("L99").Value = 3087.7144
$ws.Range("M99").Value = -4999007.5
$ws.Range("N99").Value = -6083.7144

$ws.Range("H105").Value = 1022
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 846.1111
$ws.Range("I113").Value = 802.5714
$ws.Range("K113").Value = 802.5714
$ws.Range("M113").Value = 1367.4286

$ws.Range("H126").Value = 1113625
$ws.Range("I126").Value = 5000505.5
$ws.Range("J126").Value = 3087.7144
$ws.Range("K126").Value = 15001516.5
$ws.Range("L126").Value = 9263.143199999999
$ws.Range("M126").Value = -14999046.5
$ws.Range("N126").Value = -14203.1432

$ws.Range("H132").Value = 2432.5715
$ws.Range("I132").Value = 1064.5
$ws.Range("K132").Value = 3193.5
$ws.Range("M132").Value = -663.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 829.8570999999999
$ws.Range("I5").Value = 749.5
$ws.Range("J5").Value = 862
$ws.Range("K5").Value = 2248.5
$ws.Range("L5").Value = 2586
$ws.Range("M5").Value = -2136.5
$ws.Range("N5").Value = -2810

$ws.Range("H113").Value = 138296.75
$ws.Range("I113").Value = 367001.34
$ws.Range("J113").Value = 1074
$ws.Range("K113").Value = 1101004.02
$ws.Range("L113").Value = 3222
$ws.Range("M113").Value = -1098834.02
$ws.Range("N113").Value = -7562

$ws.Range("H122").Value = 1019
$ws.Range("I122").Value = 280
$ws.Range("J122").Value = 1101.1111
$ws.Range("K122").Value = 2520
$ws.Range("L122").Value = 9909.999900000001
$ws.Range("M122").Value = -70
$ws.Range("N122").Value = -14809.9999

$ws.Range("H131").Value = 776.26
$ws.Range("J131").Value = 786.7292
$ws.Range("L131").Value = 2360.1876
$ws.Range("N131").Value = -12440.1876

$ws.Range("H135").Value = 829.8570999999999
$ws.Range("I135").Value = 749.5
$ws.Range("J135").Value = 862
$ws.Range("K135").Value = 6745.5
$ws.Range("L135").Value = 7758
$ws.Range("M135").Value = -4210.5
$ws.Range("N135").Value = -12828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1339.6
$ws.Range("I113").Value = 1099.5
$ws.Range("K113").Value = 1099.5
$ws.Range("M113").Value = 1070.5

$ws.Range("H126").Value = 2573481.2
$ws.Range("I126").Value = 2926820.5
$ws.Range("J126").Value = 335666.66
$ws.Range("K126").Value = 8780461.5
$ws.Range("L126").Value = 1006999.98
$ws.Range("M126").Value = -8777991.5
$ws.Range("N126").Value = -1011939.98

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2586.75
$ws.Range("I16").Value = 2493.7856
$ws.Range("K16").Value = 2493.7856
$ws.Range("M16").Value = -2323.7856

$ws.Range("H43").Value = 11787.111
$ws.Range("J43").Value = 11787.111
$ws.Range("L43").Value = 11787.111
$ws.Range("N43").Value = -12173.111

$ws.Range("H55").Value = 374.07693
$ws.Range("I55").Value = 296.35294
$ws.Range("J55").Value = 520.8889
$ws.Range("K55").Value = 296.35294
$ws.Range("L55").Value = 520.8889
$ws.Range("M55").Value = -123.35294
$ws.Range("N55").Value = -866.8889

$ws.Range("H68").Value = 2917.182
$ws.Range("I68").Value = 2621
$ws.Range("J68").Value = 4250
$ws.Range("K68").Value = 2621
$ws.Range("L68").Value = 4250
$ws.Range("M68").Value = -1872
$ws.Range("N68").Value = -5748

$ws.Range("H71").Value = 2917.182
$ws.Range("I71").Value = 2621
$ws.Range("J71").Value = 4250
$ws.Range("K71").Value = 13105
$ws.Range("L71").Value = 21250
$ws.Range("M71").Value = -9361
$ws.Range("N71").Value = -28738

$ws.Range("H82").Value = 1948.2354
$ws.Range("I82").Value = 1360.6666
$ws.Range("K82").Value = 1360.6666
$ws.Range("M82").Value = -999.6666

$ws.Range("H85").Value = 1948.2354
$ws.Range("I85").Value = 1360.6666
$ws.Range("K85").Value = 1360.6666
$ws.Range("M85").Value = -112.6666

$ws.Range("H122").Value = 13333.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 61799.2
$ws.Range("J108").Value = 61799.2
$ws.Range("L108").Value = 61799.2
$ws.Range("N108").Value = -69479.2

$ws.Range("H132").Value = 1615.6471
$ws.Range("I132").Value = 1264.4667
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 3793.4001
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -1263.4001
$ws.Range("N132").Value = -17808.5

$ws.Range("H136").Value = 23150898
$ws.Range("I136").Value = 39685570
$ws.Range("J136").Value = 2359.9
$ws.Range("K136").Value = 119056710
$ws.Range("L136").Value = 7079.700000000001
$ws.Range("M136").Value = -119054160
$ws.Range("N136").Value = -12179.7
